# Drop the 2009 row (row 2): every subsequent year's data shifts up by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(2).Delete()

# Append the new 2021 figures as the new last row (row 13).
$ws.Range("A13").Value = "2021年"

$vals13 = @(996.2667, 221.6175, 303.656, 658.4724, 816.9438, 1038.7434, 509.7281, 3131.7723, 11586.3388, 1602.1838, 155.9013, 137.3047, 708.21, 7113.7171, 134.8583, 9721.5681, 1141.4366, 1043.8516, 2642.7944, 1079.5435)
for ($i = 0; $i -lt $vals13.Length; $i++) {
    $col = $i + 2   # data starts at column B
    $ws.Cells.Item(13, $col).Value = $vals13[$i]
}

# Match the year-label formatting (bold, bordered, centered) used by A2:A12.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
